# Added seminar schedule H24
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Populate the "Seminar" column (F) with the new H24 seminar schedule,
# replacing the previous placeholder "-" values.
$ws.Range("F2").Value = "29.08 *08:15 - 10:00* (**LAB2**)"
$ws.Range("F3").Value = "05.09 *08:15 - 10:00* (**AUD G**)"
$ws.Range("F4").Value = "12.09 *08:15 - 10:00* (**BORCH**)"
$ws.Range("F5").Value = "19.09 *08:15 - 10:00* (**LAB2**)"
$ws.Range("F6").Value = "26.09 *08:15 - 10:00* (**LAB2**)"

# Update view state: scroll so column B is the left-most visible column,
# and move the active selection to F7.
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F7").Select()
